$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking data that must stay plain
# TEXT (it mirrors the workbook's original inlineStr cells, e.g. European
# thousands-dot notation like "62.631.47", and values such as "0.990"
# whose trailing zero would be lost if Excel auto-coerced them to a real
# number). Force text format on every Price cell we touch before writing
# the value so Excel's type-inference doesn't turn it into a Number.
$priceCells = "D2","D3","D5","D6","D8","D9","D10","D12","D14","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D29","D30","D32","D36","D37","D38","D39","D41","D43","D44","D45","D46","D47","D48","D51"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.631.47"
$ws.Range("E2").Value = "  +2.46%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.961.06"
$ws.Range("E3").Value = "  +1.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "595.77"
$ws.Range("E5").Value = "  +0.69%  "

# Row 6 - Solana
$ws.Range("D6").Value = "145.09"
$ws.Range("E6").Value = "  -0.07%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "2.957.64"
$ws.Range("E8").Value = "  +1.10%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.507"
$ws.Range("E9").Value = "  +0.34%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "7.34"
$ws.Range("E10").Value = "  +5.96%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.01%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  +1.56%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +3.61%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "33.40"
$ws.Range("E14").Value = "  -0.93%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -0.38%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "3.452.59"
$ws.Range("E16").Value = "  +1.17%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "62.568.63"
$ws.Range("E17").Value = "  +2.60%  "

# Row 18 - Polkadot
$ws.Range("D18").Value = "6.73"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "2.960.59"
$ws.Range("E19").Value = "  +1.02%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "441.70"
$ws.Range("E20").Value = "  +0.98%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").Value = "  -0.85%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "7.11"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "81.84"
$ws.Range("E24").Value = "  +0.33%  "

# Row 25 - RenderToken
$ws.Range("D25").Value = "11.13"
$ws.Range("E25").Value = "  +0.90%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "11.99"
$ws.Range("E26").Value = "  +0.93%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  -3.29%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  +0.03%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value = "2.62"
$ws.Range("E29").Value = "  +0.63%  "

# Row 30 - NEARProtocol
$ws.Range("D30").Value = "7.07"
$ws.Range("E30").Value = "  +0.87%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -6.15%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "26.58"
$ws.Range("E32").Value = "  -0.20%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -2.56%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  -0.01%  "

# Row 35 - PEPE
$ws.Range("E35").Value = "  +0.89%  "

# Row 36 - Mantle
$ws.Range("D36").Value = "0.990"
$ws.Range("E36").Value = "  -2.18%  "

# Row 37 - Filecoin
$ws.Range("D37").Value = "5.63"
$ws.Range("E37").Value = "  -0.27%  "

# Row 38 - Stacks
$ws.Range("D38").Value = "2.06"
$ws.Range("E38").Value = "  +3.26%  "

# Row 39 - OKB
$ws.Range("D39").Value = "49.63"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40 - dogwifhat
$ws.Range("E40").Value = "  -2.59%  "

# Row 41 - Cosmos
$ws.Range("D41").Value = "8.57"
$ws.Range("E41").Value = "  -0.32%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -4.21%  "

# Row 43 - TheGraph
$ws.Range("D43").Value = "0.282"
$ws.Range("E43").Value = "  -1.79%  "

# Row 44 - Arweave
$ws.Range("D44").Value = "39.09"
$ws.Range("E44").Value = "  -7.36%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.721.85"
$ws.Range("E45").Value = "  +1.18%  "

# Row 46 / 47 - VeChain and Monero swap ranking positions
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "135.43"
$ws.Range("E46").Value = "  +1.55%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0342"
$ws.Range("E47").Value = "  -1.52%  "

# Row 48 - Bittensor
$ws.Range("D48").Value = "366.27"
$ws.Range("E48").Value = "  -2.37%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  -0.19%  "

# Row 51 - InjectiveProtocol
$ws.Range("D51").Value = "23.01"
$ws.Range("E51").Value = "  -3.97%  "
